# Update cryptocurrency price/volume data as published by the GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.641.60"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").Value = "'2.296.47"
$ws.Range("E3").Value = "  +3.78%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'268.24"
$ws.Range("E5").Value = "  +1.55%  "

$ws.Range("D6").Value = "'92.59"
$ws.Range("E6").Value = "  +6.98%  "

$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  +1.73%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'0.617"
$ws.Range("E9").Value = "  +1.83%  "

$ws.Range("D10").Value = "'44.85"
$ws.Range("E10").Value = "  -3.13%  "

$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("D12").Value = "'8.01"
$ws.Range("E12").Value = "  +6.00%  "

$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").Value = "'2.635.76"
$ws.Range("E14").Value = "  +3.56%  "

$ws.Range("D15").Value = "'15.19"
$ws.Range("E15").Value = "  +3.90%  "

$ws.Range("D16").Value = "'0.846"
$ws.Range("E16").Value = "  +7.97%  "

$ws.Range("D17").Value = "'2.276.27"
$ws.Range("E17").Value = "  +4.02%  "

$ws.Range("D18").Value = "'43.632.38"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("E19").Value = "  +1.30%  "

$ws.Range("D20").Value = "'6.24"
$ws.Range("E20").Value = "  +4.45%  "

$ws.Range("D21").Value = "'71.03"
$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("D22").Value = "'2.27"
$ws.Range("E22").Value = "  -4.54%  "

$ws.Range("D23").Value = "'237.44"
$ws.Range("E23").Value = "  +2.37%  "

$ws.Range("D24").Value = "'9.66"
$ws.Range("E24").Value = "  +8.53%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").Value = "'2.50"
$ws.Range("E26").Value = "  +3.52%  "

$ws.Range("D27").Value = "'11.12"
$ws.Range("E27").Value = "  +2.67%  "

$ws.Range("D28").Value = "'3.39"
$ws.Range("E28").Value = "  -3.84%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.29"
$ws.Range("E29").Value = "  +1.20%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'39.01"
$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("D31").Value = "'22.54"
$ws.Range("E31").Value = "  +9.88%  "

$ws.Range("D32").Value = "'173.24"
$ws.Range("E32").Value = "  -1.04%  "

$ws.Range("D33").Value = "'0.0881"
$ws.Range("E33").Value = "  -1.27%  "

$ws.Range("D34").Value = "'5.52"
$ws.Range("E34").Value = "  +2.03%  "

$ws.Range("E35").Value = "  +1.64%  "

$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("D37").Value = "'4.49"
$ws.Range("E37").Value = "  +1.37%  "

$ws.Range("D38").Value = "'0.0348"
$ws.Range("E38").Value = "  -2.62%  "

$ws.Range("D39").Value = "'3.34"
$ws.Range("E39").Value = "  +2.43%  "

$ws.Range("D40").Value = "'2.33"
$ws.Range("E40").Value = "  +11.46%  "

$ws.Range("D41").Value = "'0.234"
$ws.Range("E41").Value = "  +14.89%  "

$ws.Range("D42").Value = "'12.19"
$ws.Range("E42").Value = "  -1.61%  "

$ws.Range("D43").Value = "'1.33"
$ws.Range("E43").Value = "  +18.48%  "

$ws.Range("D44").Value = "'5.43"
$ws.Range("E44").Value = "  -2.10%  "

$ws.Range("D45").Value = "'60.92"
$ws.Range("E45").Value = "  -5.65%  "

$ws.Range("D46").Value = "'8.88"
$ws.Range("E46").Value = "  +6.49%  "

$ws.Range("D47").Value = "'0.101"
$ws.Range("E47").Value = "  +2.96%  "

$ws.Range("D48").Value = "'99.65"
$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "'2.516.37"
$ws.Range("E50").Value = "  +3.59%  "

$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.430"
$ws.Range("E51").Value = "  -3.37%  "
